$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1) "Revision History" sheet: insert a new revision row (row 3) for the
#    5.0 / VistALink version-number update, pushing the existing rows down.
# ---------------------------------------------------------------------------
$rev = $wb.Worksheets.Item("Revision History")

# Insert a blank row above the current row 3 (shifts old rows 3-9 to 4-10).
$rev.Rows.Item(3).Insert()

# Copy the formatting (styles) of the row immediately below (the row that
# used to be row 3, now row 4) into the freshly-inserted row so the new row
# matches the existing "revision history" look (borders/fills/number format).
$rev.Range("A4:D4").Copy()
$rev.Range("A3:D3").PasteSpecial(-4122)  # xlPasteFormats

# Fill in the new revision entry.
$rev.Range("A3").Value = "05-20-2019"
$rev.Range("B3").Value = "5.0"
$rev.Range("C3").Value = "Updated VistALink version number "
$rev.Range("D3").Value = "Donald Fowlds"

# Leave the cursor on the newly added row, as in the authored edit.
$rev.Range("A3").Select()

# ---------------------------------------------------------------------------
# 2) "Cover" sheet: bump the footer version/date, and update the VistALink
#    component build timestamp inside the "Build 10 includes..." paragraph.
# ---------------------------------------------------------------------------
$cover = $wb.Worksheets.Item("Cover")

# Update the VistALink build id timestamp (084015 -> 085649) in the
# multi-run "Build 10 includes 4 components" cell, re-applying the run-level
# formatting (bold / size / font) that the cell already used.
$buildCell = $cover.Range("A3")
$text = $buildCell.Text
$newText = $text.Replace("084015", "085649")
$buildCell.Value = $newText

$run1 = $buildCell.Characters(1, 11)
$run1.Font.Bold = $true
$run1.Font.Size = 12
$run1.Font.Name = "Arial"

$run2 = $buildCell.Characters(12, 1)
$run2.Font.Bold = $true
$run2.Font.Size = 12
$run2.Font.Name = "Calibri"

$run3 = $buildCell.Characters(13, 41)
$run3.Font.Bold = $true
$run3.Font.Size = 12
$run3.Font.Name = "Arial"

$run4 = $buildCell.Characters(54, 1)
$run4.Font.Bold = $true
$run4.Font.Size = 12
$run4.Font.Name = "Calibri"

$run5 = $buildCell.Characters(55, 43)
$run5.Font.Bold = $true
$run5.Font.Size = 12
$run5.Font.Name = "Arial"

$run6 = $buildCell.Characters(98, 1)
$run6.Font.Bold = $true
$run6.Font.Size = 12
$run6.Font.Name = "Calibri"

$run7 = $buildCell.Characters(99, 54)
$run7.Font.Bold = $true
$run7.Font.Size = 12
$run7.Font.Name = "Arial"

$run8 = $buildCell.Characters(153, 1)
$run8.Font.Bold = $true
$run8.Font.Size = 12
$run8.Font.Name = "Calibri"

$run9 = $buildCell.Characters(154, 36)
$run9.Font.Bold = $true
$run9.Font.Size = 12
$run9.Font.Name = "Arial"

# Update the footer date/version lines.
$cover.Range("A15").Value = "May 20, 2019"
$cover.Range("A16").Value = "Version 5.0"

# Leave the cursor on A2, as in the authored edit, and make sure "Cover"
# is the sheet that is active/selected when the workbook is saved.
$cover.Range("A2").Select()
$cover.Activate()
